# 1.0 parse comm equipment and save to excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# F1 (new "Дата" column) should inherit the bordered/bold/centered + date-number-format
# style that B1 ("data") used to carry - grab it before B1's own style is changed below.
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B1 no longer holds dates - give it the plain header style shared by the other labels.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row values: translate / relabel columns ---
$ws.Range("A1").Value = "Модель"
$ws.Range("B1").Value = "Москва"
$ws.Range("C1").Value = "Санкт-Петербург"
$ws.Range("D1").Value = "Новосибирск"
$ws.Range("E1").Value = "В транзите"
$ws.Range("F1").Value = "Дата"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 34.43
$ws.Columns.Item(2).ColumnWidth = 16.86
$ws.Columns.Item(3).ColumnWidth = 16.86
$ws.Columns.Item(4).ColumnWidth = 16.86
$ws.Columns.Item(5).ColumnWidth = 10.71
$ws.Columns.Item(6).ColumnWidth = 16.14

# --- Selection as left by the author ---
$ws.Range("D3").Select()
